$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C ("Quantite") shifting the old "Quantité Ingrédient"
# column (and everything after it) one column to the right.
$ws.Columns("C").Insert()

# New header for the inserted column.
$ws.Range("C1").Value = "Quantite"

# Fill the new numeric "Quantite" column with the quantity extracted from the
# existing (now shifted to column D) quantity-with-unit strings.
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 125
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 2
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 2
$ws.Range("C8").Value = 235
$ws.Range("C9").Value = 100
$ws.Range("C10").Value = 5
$ws.Range("C11").Value = 9
$ws.Range("C12").Value = 10
$ws.Range("C13").Value = 0
$ws.Range("C14").Value = 2
$ws.Range("C15").Value = 2

# Restore the current selection to match the edited location.
$ws.Range("C15").Select()
